# Update the interactive PEBCOM map: the record that was in row 59
# (Caso 5940 - SANCHEZ DE LORIA 1406 / OT 807044148) was removed from the
# source feed. Delete that entire worksheet row so every following record
# shifts up by one position (row 60 -> 59, row 61 -> 60, ... row 81 -> 80),
# shrinking the used range from A1:P81 down to A1:P80.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEBCOM")

$ws.Rows.Item(59).Delete()
